# Applies the "dodelani nove truktury ip setting" commit to the workbook.
$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet: ip_address_list
# -------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")

# Row3 / D3: drop the stray trailing "j" on the separator line
$ws1.Range("D3").Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.205.267`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"

# Row7 / D7: fix the camera IP (186 -> 18)
$ws1.Range("D7").Value = "Kamera VS-S160MX :192.168.0.18"

# Row8 / D8: remove the stray note entirely
$ws1.Range("D8").ClearContents()

# Row11 / D11: drop the stray trailing "f"
$ws1.Range("D11").Value = "XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28.117"

# Row12 / D12: trim the NAS note
$ws1.Range("D12").Value = "NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 `nUser:jhvadmin Pass`n123TPV456"

# -------------------------------------------------------------------
# Sheet: ip_address_fav_list
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ip_address_fav_list")

# Row4 / D4: same "OP:" fix as above
$ws2.Range("D4").Value = "XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28.117"

# Row5 / D5: same NAS-note trim as above
$ws2.Range("D5").Value = "NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 `nUser:jhvadmin Pass`n123TPV456"

# -------------------------------------------------------------------
# Sheet: disk_list - restructured (rows reordered / merged / trimmed,
# one row dropped -> dimension goes from A1:F8 to A1:F7)
# -------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("disk_list")
$ws3.Range("A1:F8").ClearContents()

$ws3.Range("A1").Value = "514_Teleflex"
$ws3.Range("B1").Value = "T"
$ws3.Range("C1").Value = "\\192.168.14.245\Data\Kamery"
$ws3.Range("D1").Value = "Vision"
$ws3.Range("E1").Value = "*Jhv2708"

$ws3.Range("A2").Value = "515_ZF"
$ws3.Range("B2").Value = "Z"
$ws3.Range("C2").Value = "\\10.9.250.100\08_Project_ZF_515\kamery"
$ws3.Range("D2").Value = "jhvadmin"
$ws3.Range("E2").Value = "jhvadm1n"

$ws3.Range("A3").Value = "Domaci Nas"
$ws3.Range("B3").Value = "S"
$ws3.Range("C3").Value = "\\192.168.1.20\Data"

$ws3.Range("A4").Value = "518_Valeo II"
$ws3.Range("B4").Value = "V"
$ws3.Range("C4").Value = "\\192.168.1.10\10_vision"
$ws3.Range("D4").Value = "jhv_vision"
$ws3.Range("E4").Value = "Jhv*2708"
$ws3.Range("F4").Value = "Druha sít, ixon"

$ws3.Range("A5").Value = "518_Valeo"
$ws3.Range("B5").Value = "V"
$ws3.Range("C5").Value = "\\192.168.208.200\10_vision"
$ws3.Range("D5").Value = "jhv_vision"
$ws3.Range("E5").Value = "Jhv*2708"
$ws3.Range("F5").Value = "první sít, ixon`n\\192.168.208.200\10_vision"

$ws3.Range("A6").Value = "474_B Austin"
$ws3.Range("B6").Value = "P"
$ws3.Range("C6").Value = "\\10.96.205.166\DATA"
$ws3.Range("D6").Value = "jhv_vision"
$ws3.Range("E6").Value = "*Jhv2708"
$ws3.Range("F6").Value = "10.96.205.166`nVisionNas_474B`t`n`t`t`t`t`t`tuser:JHV_Vision, omron `nPass:*Jhv2708"

$ws3.Range("A7").Value = "xfdx"
$ws3.Range("B7").Value = "P"
$ws3.Range("C7").Value = "\\192.168.000.000\"
$ws3.Range("D7").Value = "ss"

# -------------------------------------------------------------------
# Sheet: projects_bin2 (hidden) - now fully emptied
# -------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("projects_bin2")
$ws5.Range("A3:F4").ClearContents()

# -------------------------------------------------------------------
# Sheet: Settings_recources
# -------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Settings_recources")

# B3 must stay a literal text "False" (not a real boolean)
$ws6.Range("B3").Value = "'False"

# B20 must stay literal text "85" (not a real number)
$ws6.Range("B20").Value = "'85"

$ws6.Range("B28").Value = "xlsx"

$ws6.Range("A30").Value = "nastavení zoomu celé aplikace (default: 100 %)"
$ws6.Range("B30").Value = 100
